# Auto-generated Excel COM-interop script
# Updates currentAveragePrice / Leve price & profit columns (H-N)
# across several worksheets, per scheduled-runner data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 216.2
$ws.Cells.Item(38, 10).Value = 38
$ws.Cells.Item(38, 12).Value = 114
$ws.Cells.Item(38, 14).Value = -858

$ws.Cells.Item(41, 8).Value = 2255.7778
$ws.Cells.Item(41, 9).Value = 3117.1667
$ws.Cells.Item(41, 11).Value = 3117.1667
$ws.Cells.Item(41, 13).Value = -2677.1667

$ws.Cells.Item(113, 8).Value = 21663.6
$ws.Cells.Item(113, 9).Value = 7159.25
$ws.Cells.Item(113, 10).Value = 31333.166
$ws.Cells.Item(113, 11).Value = 7159.25
$ws.Cells.Item(113, 12).Value = 31333.166
$ws.Cells.Item(113, 13).Value = -3905.25
$ws.Cells.Item(113, 14).Value = -37841.166

$ws.Cells.Item(129, 8).Value = 2605081.2
$ws.Cells.Item(129, 9).Value = 695.0769
$ws.Cells.Item(129, 11).Value = 2085.2307
$ws.Cells.Item(129, 13).Value = 2914.7693

$ws.Cells.Item(132, 8).Value = 2059.04
$ws.Cells.Item(132, 9).Value = 1917.9048
$ws.Cells.Item(132, 11).Value = 5753.7144
$ws.Cells.Item(132, 13).Value = -3223.7144

$ws.Cells.Item(138, 8).Value = 4197.6
$ws.Cells.Item(138, 9).Value = 2564.5715
$ws.Cells.Item(138, 10).Value = 5076.923
$ws.Cells.Item(138, 11).Value = 7693.7145
$ws.Cells.Item(138, 12).Value = 15230.769
$ws.Cells.Item(138, 13).Value = -2553.7145
$ws.Cells.Item(138, 14).Value = -25510.769

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(33, 8).Value = 4866.3335
$ws.Cells.Item(33, 9).Value = 4866.3335
$ws.Cells.Item(33, 11).Value = 4866.3335
$ws.Cells.Item(33, 13).Value = -4537.3335

$ws.Cells.Item(63, 8).Value = 717043.4399999999
$ws.Cells.Item(63, 9).Value = 1800.9
$ws.Cells.Item(63, 10).Value = 2505149.8
$ws.Cells.Item(63, 11).Value = 1800.9
$ws.Cells.Item(63, 12).Value = 2505149.8
$ws.Cells.Item(63, 13).Value = -1114.9
$ws.Cells.Item(63, 14).Value = -2506521.8

$ws.Cells.Item(66, 8).Value = 717043.4399999999
$ws.Cells.Item(66, 9).Value = 1800.9
$ws.Cells.Item(66, 10).Value = 2505149.8
$ws.Cells.Item(66, 11).Value = 9004.5
$ws.Cells.Item(66, 12).Value = 12525749
$ws.Cells.Item(66, 13).Value = -5572.5
$ws.Cells.Item(66, 14).Value = -12532613

$ws.Cells.Item(110, 8).Value = 8388.888999999999
$ws.Cells.Item(110, 9).Value = 7071.4287
$ws.Cells.Item(110, 11).Value = 7071.4287
$ws.Cells.Item(110, 13).Value = -5026.4287

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 4360
$ws.Cells.Item(134, 9).Value = 4419.6597
$ws.Cells.Item(134, 11).Value = 13258.9791
$ws.Cells.Item(134, 13).Value = -10723.9791

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 4772.923
$ws.Cells.Item(16, 9).Value = 3934.9
$ws.Cells.Item(16, 11).Value = 3934.9
$ws.Cells.Item(16, 13).Value = -3647.9

$ws.Cells.Item(31, 8).Value = 5229.3335
$ws.Cells.Item(31, 9).Value = 3746.1
$ws.Cells.Item(31, 10).Value = 6101.8237
$ws.Cells.Item(31, 11).Value = 3746.1
$ws.Cells.Item(31, 12).Value = 6101.8237
$ws.Cells.Item(31, 13).Value = -3451.1
$ws.Cells.Item(31, 14).Value = -6691.8237

$ws.Cells.Item(34, 8).Value = 5229.3335
$ws.Cells.Item(34, 9).Value = 3746.1
$ws.Cells.Item(34, 10).Value = 6101.8237
$ws.Cells.Item(34, 11).Value = 3746.1
$ws.Cells.Item(34, 12).Value = 6101.8237
$ws.Cells.Item(34, 13).Value = -3544.1
$ws.Cells.Item(34, 14).Value = -6505.8237

$ws.Cells.Item(35, 8).Value = 12105
$ws.Cells.Item(35, 9).Value = 2512.5
$ws.Cells.Item(35, 11).Value = 2512.5
$ws.Cells.Item(35, 13).Value = -2218.5

$ws.Cells.Item(99, 8).Value = 7139.8
$ws.Cells.Item(99, 10).Value = 7000
$ws.Cells.Item(99, 12).Value = 7000
$ws.Cells.Item(99, 14).Value = -9996

$ws.Cells.Item(113, 8).Value = 4772.923
$ws.Cells.Item(113, 9).Value = 3934.9
$ws.Cells.Item(113, 11).Value = 3934.9
$ws.Cells.Item(113, 13).Value = -1764.9

$ws.Cells.Item(126, 8).Value = 7139.8
$ws.Cells.Item(126, 10).Value = 7000
$ws.Cells.Item(126, 12).Value = 21000
$ws.Cells.Item(126, 14).Value = -25940

$ws.Cells.Item(134, 8).Value = 7763.5454
$ws.Cells.Item(134, 9).Value = 5759.8
$ws.Cells.Item(134, 11).Value = 17279.4
$ws.Cells.Item(134, 13).Value = -14744.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(23, 8).Value = 292.33334
$ws.Cells.Item(23, 9).Value = 344
$ws.Cells.Item(23, 10).Value = 266.5
$ws.Cells.Item(23, 11).Value = 1032
$ws.Cells.Item(23, 12).Value = 799.5
$ws.Cells.Item(23, 13).Value = -797
$ws.Cells.Item(23, 14).Value = -1269.5

$ws.Cells.Item(70, 8).Value = 800
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 800
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 2400
$ws.Cells.Item(70, 13).Value = $null
$ws.Cells.Item(70, 14).Value = -3030

$ws.Cells.Item(73, 8).Value = 800
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 800
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 12).Value = 2400
$ws.Cells.Item(73, 13).Value = $null
$ws.Cells.Item(73, 14).Value = -4584

$ws.Cells.Item(132, 8).Value = 1781.4375
$ws.Cells.Item(132, 10).Value = 1921.8462
$ws.Cells.Item(132, 12).Value = 17296.6158
$ws.Cells.Item(132, 14).Value = -22356.6158

$ws.Cells.Item(134, 8).Value = 971.63635
$ws.Cells.Item(134, 9).Value = 922.6667
$ws.Cells.Item(134, 11).Value = 2768.0001
$ws.Cells.Item(134, 13).Value = 2301.9999

$ws.Cells.Item(137, 8).Value = 3980.2
$ws.Cells.Item(137, 9).Value = 3609.7778
$ws.Cells.Item(137, 10).Value = 4283.273
$ws.Cells.Item(137, 11).Value = 10829.3334
$ws.Cells.Item(137, 12).Value = 12849.819
$ws.Cells.Item(137, 13).Value = -5729.3334
$ws.Cells.Item(137, 14).Value = -23049.819

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(3, 8).Value = 2852.5
$ws.Cells.Item(3, 9).Value = 700
$ws.Cells.Item(3, 10).Value = 5005
$ws.Cells.Item(3, 11).Value = 700
$ws.Cells.Item(3, 12).Value = 5005
$ws.Cells.Item(3, 13).Value = -584
$ws.Cells.Item(3, 14).Value = -5237

$ws.Cells.Item(7, 8).Value = 104289.9
$ws.Cells.Item(7, 10).Value = 104289.9
$ws.Cells.Item(7, 12).Value = 104289.9
$ws.Cells.Item(7, 14).Value = -104513.9

$ws.Cells.Item(8, 8).Value = 104289.9
$ws.Cells.Item(8, 10).Value = 104289.9
$ws.Cells.Item(8, 12).Value = 104289.9
$ws.Cells.Item(8, 14).Value = -104567.9

$ws.Cells.Item(14, 8).Value = 4600.4287
$ws.Cells.Item(14, 9).Value = 2466
$ws.Cells.Item(14, 10).Value = 6201.25
$ws.Cells.Item(14, 11).Value = 2466
$ws.Cells.Item(14, 12).Value = 6201.25
$ws.Cells.Item(14, 13).Value = -2298
$ws.Cells.Item(14, 14).Value = -6537.25

$ws.Cells.Item(102, 8).Value = 2879.3872
$ws.Cells.Item(102, 9).Value = 1509.5769
$ws.Cells.Item(102, 10).Value = 10002.4
$ws.Cells.Item(102, 11).Value = 1509.5769
$ws.Cells.Item(102, 12).Value = 10002.4
$ws.Cells.Item(102, 13).Value = 112.4231
$ws.Cells.Item(102, 14).Value = -13246.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(42, 8).Value = 12500
$ws.Cells.Item(42, 10).Value = 20000
$ws.Cells.Item(42, 12).Value = 20000
$ws.Cells.Item(42, 14).Value = -21126

$ws.Cells.Item(49, 8).Value = 12500
$ws.Cells.Item(49, 10).Value = 20000
$ws.Cells.Item(49, 12).Value = 20000
$ws.Cells.Item(49, 14).Value = -20294

$ws.Cells.Item(61, 8).Value = 37444.715
$ws.Cells.Item(61, 9).Value = 40245.31
$ws.Cells.Item(61, 10).Value = 1037
$ws.Cells.Item(61, 11).Value = 40245.31
$ws.Cells.Item(61, 12).Value = 1037
$ws.Cells.Item(61, 13).Value = -40043.31
$ws.Cells.Item(61, 14).Value = -1441

$ws.Cells.Item(68, 8).Value = 5258.1816
$ws.Cells.Item(68, 10).Value = 8877.799999999999
$ws.Cells.Item(68, 12).Value = 8877.799999999999
$ws.Cells.Item(68, 14).Value = -10375.8

$ws.Cells.Item(71, 8).Value = 5258.1816
$ws.Cells.Item(71, 10).Value = 8877.799999999999
$ws.Cells.Item(71, 12).Value = 44389
$ws.Cells.Item(71, 14).Value = -51877

$ws.Cells.Item(76, 8).Value = 45000
$ws.Cells.Item(76, 10).Value = 45000
$ws.Cells.Item(76, 12).Value = 45000
$ws.Cells.Item(76, 14).Value = -45676

$ws.Cells.Item(79, 8).Value = 45000
$ws.Cells.Item(79, 10).Value = 45000
$ws.Cells.Item(79, 12).Value = 45000
$ws.Cells.Item(79, 14).Value = -47340

$ws.Cells.Item(100, 8).Value = 4006.0952
$ws.Cells.Item(100, 9).Value = 2193.4546
$ws.Cells.Item(100, 11).Value = 2193.4546
$ws.Cells.Item(100, 13).Value = -1652.4546

$ws.Cells.Item(113, 8).Value = 37444.715
$ws.Cells.Item(113, 9).Value = 40245.31
$ws.Cells.Item(113, 10).Value = 1037
$ws.Cells.Item(113, 11).Value = 40245.31
$ws.Cells.Item(113, 12).Value = 1037
$ws.Cells.Item(113, 13).Value = -38075.31
$ws.Cells.Item(113, 14).Value = -5377

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(15, 8).Value = 7494
$ws.Cells.Item(15, 10).Value = 7494
$ws.Cells.Item(15, 12).Value = 7494
$ws.Cells.Item(15, 14).Value = -8070

$ws.Cells.Item(61, 8).Value = 35250
$ws.Cells.Item(61, 9).Value = 30500
$ws.Cells.Item(61, 10).Value = 40000
$ws.Cells.Item(61, 11).Value = 30500
$ws.Cells.Item(61, 12).Value = 40000
$ws.Cells.Item(61, 13).Value = -30208
$ws.Cells.Item(61, 14).Value = -40584

$ws.Cells.Item(74, 8).Value = 20899.334
$ws.Cells.Item(74, 9).Value = 12000
$ws.Cells.Item(74, 10).Value = 22679.2
$ws.Cells.Item(74, 11).Value = 12000
$ws.Cells.Item(74, 12).Value = 22679.2
$ws.Cells.Item(74, 13).Value = -11064
$ws.Cells.Item(74, 14).Value = -24551.2

$ws.Cells.Item(77, 8).Value = 20899.334
$ws.Cells.Item(77, 9).Value = 12000
$ws.Cells.Item(77, 10).Value = 22679.2
$ws.Cells.Item(77, 11).Value = 36000
$ws.Cells.Item(77, 12).Value = 68037.60000000001
$ws.Cells.Item(77, 13).Value = -31320
$ws.Cells.Item(77, 14).Value = -77397.60000000001

$ws.Cells.Item(106, 8).Value = 0
$ws.Cells.Item(106, 10).Value = 0
$ws.Cells.Item(106, 12).Value = 0
$ws.Cells.Item(106, 14).Value = $null

$ws.Cells.Item(109, 8).Value = 41999.75
$ws.Cells.Item(109, 10).Value = 41999.75
$ws.Cells.Item(109, 12).Value = 41999.75
$ws.Cells.Item(109, 14).Value = -44773.75

$ws.Cells.Item(113, 8).Value = 745.625
$ws.Cells.Item(113, 9).Value = 620.8
$ws.Cells.Item(113, 11).Value = 1862.4
$ws.Cells.Item(113, 13).Value = 307.6000000000001

$ws.Cells.Item(122, 8).Value = 3583.3684
$ws.Cells.Item(122, 9).Value = 2049.1667
$ws.Cells.Item(122, 10).Value = 6213.4287
$ws.Cells.Item(122, 11).Value = 6147.500100000001
$ws.Cells.Item(122, 12).Value = 18640.2861
$ws.Cells.Item(122, 13).Value = -3697.500100000001
$ws.Cells.Item(122, 14).Value = -23540.2861

$ws.Cells.Item(132, 8).Value = 5629.278
$ws.Cells.Item(132, 9).Value = 5080.4243
$ws.Cells.Item(132, 11).Value = 15241.2729
$ws.Cells.Item(132, 13).Value = -12711.2729

$ws.Cells.Item(136, 8).Value = 5408.905
$ws.Cells.Item(136, 9).Value = 4921.6113
$ws.Cells.Item(136, 11).Value = 14764.8339
$ws.Cells.Item(136, 13).Value = -12214.8339
